$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------------
# 1) Capture the ORIGINAL special "last row" formatting (row 149) before
#    we touch it, so it can be re-applied to the new last row (152).
# ---------------------------------------------------------------------
$ws.Range("G152").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("A149:K149").Copy()
$ws.Range("A152:K152").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Row 149 becomes a normal data row (same look as row 148).
# ---------------------------------------------------------------------
$ws.Range("A148:K148").Copy()
$ws.Range("A149:K149").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) New rows 150 and 151 are normal data rows too.
# ---------------------------------------------------------------------
$ws.Range("G150").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("A148:K148").Copy()
$ws.Range("A150:K150").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("G151").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("A148:K148").Copy()
$ws.Range("A151:K151").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4) Grow Table1 to cover the three new rows.
# ---------------------------------------------------------------------
$lo.Resize($ws.Range("A8:K152"))

# ---------------------------------------------------------------------
# 5) Cell value edits inside the table body.
#    (Shared strings must be entered in this order so new entries land
#    at indices 107, 108, 109 - SP(2-0-0), 10/5,6/2023, 2024.)
# ---------------------------------------------------------------------
$ws.Range("C119").Value = 1.25

$ws.Range("A122").Value = 45200
$ws.Range("B122").Value = "SP(2-0-0)"
$ws.Range("C122").Value = 1.25
$ws.Range("K122").Value = "10/5,6/2023"

$ws.Range("B123").Value = "SL(1-0-0)"
$ws.Range("H123").Value = 1
$ws.Range("K119").Copy()
$ws.Range("K123").PasteSpecial(-4122)   # xlPasteFormats (reuse date style)
$excel.CutCopyMode = 0
$ws.Range("K123").Value = 45210

$ws.Range("B124").Value = "SL(1-0-0)"
$ws.Range("H124").Value = 1
$ws.Range("K119").Copy()
$ws.Range("K124").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("K124").Value = 45223

$ws.Range("B125").Value = "VL(1-0-0)"
$ws.Range("D125").Value = 1
$ws.Range("K119").Copy()
$ws.Range("K125").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("K125").Value = 45225

$ws.Range("A126").Value = 45231
$ws.Range("A127").Value = 45261

$ws.Range("A103").Copy()
$ws.Range("A128").PasteSpecial(-4122)   # reuse "year" style (quote-prefixed text)
$excel.CutCopyMode = 0
$ws.Range("A128").Value = "'2024"

$ws.Range("A129").Value = 45292
$ws.Range("A130").Value = 45323
$ws.Range("A131").Value = 45352
$ws.Range("A132").Value = 45383
$ws.Range("A133").Value = 45413
$ws.Range("A134").Value = 45444
$ws.Range("A135").Value = 45474
$ws.Range("A136").Value = 45505
$ws.Range("A137").Value = 45536
$ws.Range("A138").Value = 45566
$ws.Range("A139").Value = 45597
$ws.Range("A140").Value = 45627

# ---------------------------------------------------------------------
# 6) Recalculate so the BALANCE formulas in E9/I9 (and the G-column
#    helper cells) pick up the new EARNED/Absence figures.
# ---------------------------------------------------------------------
$excel.CalculateFullRebuild()

# ---------------------------------------------------------------------
# 7) Restore the on-screen selection to K125 (bottom pane), as recorded
#    in the saved workbook view.
# ---------------------------------------------------------------------
$ws.Range("K125").Select()
